$wb = $excel.ActiveWorkbook

# --- Sheet: TwoLineHeader (sheet2) ---
# Add a new row of data (row 5) below the existing data
$ws2 = $wb.Worksheets.Item("TwoLineHeader")

$ws2.Range("B5").Value = 11
$ws2.Range("C5").Value = 22
$ws2.Range("D5").Value = 33
$ws2.Range("E5").Value = 444

# Match number formatting (style index 1) used by columns B and D in row 4
$ws2.Range("B5").NumberFormat = "0"
$ws2.Range("D5").NumberFormat = "0"

# Move the active cell / selection to E5 on this sheet
$ws2.Range("E5").Select()

# --- Sheet: dataSheet (sheet1) ---
# Move the active cell / selection to E20 on this sheet
$ws1 = $wb.Worksheets.Item("dataSheet")
$ws1.Range("E20").Select()

# Re-activate the TwoLineHeader sheet as the last active tab
$ws2.Activate()
